$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (JA1): Run batch filled in
$ws.Range("C3").Value = "UWPR Apr 2017"

# Row 4 (JA2): sample name fixed from "JA2:" to "JA2", run batch filled in,
# spectra counts updated, % DNO peptides modified calculated
$ws.Range("B4").Value = "JA2"
$ws.Range("C4").Value = "UWPR Apr 2017"
$ws.Range("J4").Value = 2462
$ws.Range("K4").Value = 1467
$ws.Range("L4").Value = 3437
$ws.Range("M4").Value = 0.114906832298136

# Row 5 (JA3): Run batch filled in
$ws.Range("C5").Value = "UWPR Apr 2017"

# Row 6 (JA4): sample name fixed from "JA4:" to "JA4", run batch filled in,
# spectra counts updated, % DNO peptides modified calculated
$ws.Range("B6").Value = "JA4"
$ws.Range("C6").Value = "UWPR Apr 2017"
$ws.Range("J6").Value = 1862
$ws.Range("K6").Value = 1243
$ws.Range("L6").Value = 2854
$ws.Range("M6").Value = 0.094059405940594004

# Row 7 (JA5): Run batch filled in
$ws.Range("C7").Value = "UWPR Apr 2017"

# Row 8 (JA6): Run batch filled in
$ws.Range("C8").Value = "UWPR Apr 2017"

# Window / view changes
$excel.ActiveWindow.TopLeftCell = $ws.Range("C1")
$ws.Range("M7").Select()
